# Refactor: rename the sheet, add a named range over a small lookup
# table, and populate the new lookup rows (J16:K18).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet was a "- Copy" duplicate; rename back to the canonical name.
$ws.Name = "test"

# New lookup table used by the named range below.
$ws.Range("J16").Value = "nam_ran_col1"
$ws.Range("K16").Value = "nam_ran_col1"
$ws.Range("J17").Value = "val1"
$ws.Range("J18").Value = "val2"
$ws.Range("K17").Value = "val3"
$ws.Range("K18").Value = "val4"

# Named range covering the new lookup table.
$wb.Names.Add("xlsx_named_range1", "=test!`$J`$16:`$K`$18")

# Leave the selection where the user last clicked while doing this.
$ws.Range("X15").Select() | Out-Null
